# Auto-applies the cryptos.xlsx price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.794.90"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "'3.104.19"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'525.69"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "'141.68"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'3.106.06"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").Value = "'0.437"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'7.27"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'0.382"
$ws.Range("E12").Value = "  +3.54%  "
$ws.Range("D13").Value = "'3.638.48"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").Value = "'26.26"
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "'57.898.83"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "'3.101.77"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "'6.07"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "'12.82"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'8.06"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'337.60"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'0.505"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'66.23"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'0.0₃0906"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "'6.59"
$ws.Range("E29").Value = "  +4.21%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'7.22"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'20.98"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.20"
$ws.Range("E34").Value = "  +3.60%  "
$ws.Range("D35").Value = "'154.16"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("D37").Value = "'6.03"
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("D38").Value = "'27.08"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").Value = "'0.0663"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").Value = "'3.150.31"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("D43").Value = "'3.88"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "'36.78"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  +6.75%  "
$ws.Range("D47").Value = "'2.288.92"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").Value = "'0.0256"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'20.65"
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'0.962"
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("D51").Value = "'6.01"
$ws.Range("E51").Value = "  +2.61%  "
